# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "VALOR MORA" total and "Cant. Periodos" count
$ws.Range("E11").Value = 877691
$ws.Range("F13").Value = 12

# 2. Insert one extra data row so the table grows from 11 to 12 period rows.
#    Insert before the current last row (row 26) so that row shifts down to 27
#    and keeps its special "bottom border" formatting; the newly inserted
#    row 26 picks up the regular interior-row formatting (copied from row 25).
$ws.Rows("26").Insert()
$ws.Range("B25:J25").Copy()
$ws.Range("B26:J26").PasteSpecial(-4122)

# 3. Re-write the whole worker table (rows 16-27) with the updated/sorted data:
#    JHON JAIRO's single period moves to the top (row 16), and JOSE MARIA's
#    periods are now listed in ascending order (2410..2508), including the
#    newly added period 2508 in the last row (27).
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1042428519"
$ws.Range("D16").Value = "JHON JAIRO SARMIENTO BOLAÃ?O"
$ws.Range("E16").Value = "2402"
$ws.Range("F16").Value = 2333
$ws.Range("G16").Value = 1750000

$periods = @("2410","2411","2412","2501","2502","2503","2504","2505","2506","2507","2508")
$row = 17
foreach ($p in $periods) {
    $ws.Range("B$row").Value = "CC"
    $ws.Range("C$row").Value = "1082932914"
    $ws.Range("D$row").Value = "JOSE MARIA AREVALO CARRILLO"
    $ws.Range("E$row").Value = $p
    $ws.Range("F$row").Value = 79578
    $ws.Range("G$row").Value = 1989456
    $row = $row + 1
}

Write-Output "edit complete"
